# Auto-generated script applying cryptos.xlsx data refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "26.464.54"
$ws.Cells.Item(2, 5).Value = "  +2.69%  "
$ws.Cells.Item(3, 4).Value = "1.728.82"
$ws.Cells.Item(3, 5).Value = "  +3.12%  "
$ws.Cells.Item(4, 4).NumberFormat = "@"
$ws.Cells.Item(4, 4).Value = "0.9997"
$ws.Cells.Item(4, 5).Value = "  +0.11%  "
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "243.71"
$ws.Cells.Item(5, 5).Value = "  +2.65%  "
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "1.000"
$ws.Cells.Item(6, 5).Value = "  +0.03%  "
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = "0.4798"
$ws.Cells.Item(7, 5).Value = "  +3.66%  "
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = "0.2668"
$ws.Cells.Item(8, 5).Value = "  +2.79%  "
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = "0.06226"
$ws.Cells.Item(9, 5).Value = "  +1.36%  "
$ws.Cells.Item(10, 4).Value = "1.736.57"
$ws.Cells.Item(10, 5).Value = "  +3.61%  "
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = "0.07115"
$ws.Cells.Item(11, 5).Value = "  +1.71%  "
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = "15.72"
$ws.Cells.Item(12, 5).Value = "  +5.43%  "
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = "0.6173"
$ws.Cells.Item(13, 5).Value = "  +7.05%  "
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = "4.542"
$ws.Cells.Item(14, 5).Value = "  +4.08%  "
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = "77.09"
$ws.Cells.Item(15, 5).Value = "  +2.24%  "
$ws.Cells.Item(16, 5).Value = "  +0.00%  "
$ws.Cells.Item(17, 4).Value = "26.486.21"
$ws.Cells.Item(17, 5).Value = "  +2.80%  "
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = "0.9998"
$ws.Cells.Item(18, 5).Value = "  +0.00%  "
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = "0.000006929"
$ws.Cells.Item(19, 5).Value = "  +3.45%  "
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = "11.72"
$ws.Cells.Item(20, 5).Value = "  +2.53%  "
$ws.Cells.Item(21, 4).Value = "1.955.36"
$ws.Cells.Item(21, 5).Value = "  +3.67%  "
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "4.565"
$ws.Cells.Item(22, 5).Value = "  +2.19%  "
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "8.916"
$ws.Cells.Item(23, 5).Value = "  +2.86%  "
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = "5.317"
$ws.Cells.Item(24, 5).Value = "  +1.72%  "
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = "136.41"
$ws.Cells.Item(25, 5).Value = "  +1.69%  "
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = "15.33"
$ws.Cells.Item(26, 5).Value = "  +2.52%  "
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = "1.791"
$ws.Cells.Item(27, 5).Value = "  +4.31%  "
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = "1.411"
$ws.Cells.Item(28, 5).Value = "  +1.52%  "
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = "106.68"
$ws.Cells.Item(29, 5).Value = "  +2.08%  "
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = "3.987"
$ws.Cells.Item(30, 5).Value = "  +1.01%  "
$ws.Cells.Item(31, 5).Value = "  +3.86%  "
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = "3.739"
$ws.Cells.Item(32, 5).Value = "  +3.39%  "
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = "0.04561"
$ws.Cells.Item(33, 5).Value = "  +5.00%  "
$ws.Cells.Item(34, 2).Value = "Frax"
$ws.Cells.Item(34, 3).Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = "0.9997"
$ws.Cells.Item(34, 5).Value = "  +0.11%  "
$ws.Cells.Item(35, 2).Value = "HuobiToken"
$ws.Cells.Item(35, 3).Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = "2.614"
$ws.Cells.Item(35, 5).Value = "  +0.60%  "
$ws.Cells.Item(36, 2).Value = "ImmutableX"
$ws.Cells.Item(36, 3).Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = "0.6411"
$ws.Cells.Item(36, 5).Value = "  +5.46%  "
$ws.Cells.Item(37, 2).Value = "ARBITRUM"
$ws.Cells.Item(37, 3).Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = "0.9907"
$ws.Cells.Item(37, 5).Value = "  +4.18%  "
$ws.Cells.Item(38, 2).Value = "TrustWalletToken"
$ws.Cells.Item(38, 3).Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = "0.9422"
$ws.Cells.Item(38, 5).Value = "  +0.97%  "
$ws.Cells.Item(39, 2).Value = "RenderToken"
$ws.Cells.Item(39, 3).Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = "1.993"
$ws.Cells.Item(39, 5).Value = "  +7.16%  "
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = "2.414"
$ws.Cells.Item(40, 5).Value = "  -1.27%  "
$ws.Cells.Item(41, 2).Value = "Quant"
$ws.Cells.Item(41, 3).Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = "107.11"
$ws.Cells.Item(41, 5).Value = "  -1.77%  "
$ws.Cells.Item(42, 2).Value = "PaxDollar"
$ws.Cells.Item(42, 3).Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = "1.006"
$ws.Cells.Item(42, 5).Value = "  +0.75%  "
$ws.Cells.Item(43, 2).Value = "VeChain"
$ws.Cells.Item(43, 3).Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = "0.01501"
$ws.Cells.Item(43, 5).Value = "  +3.52%  "
$ws.Cells.Item(44, 2).Value = "FraxShare"
$ws.Cells.Item(44, 3).Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = "5.644"
$ws.Cells.Item(44, 5).Value = "  +11.87%  "
$ws.Cells.Item(45, 2).Value = "TheSandbox"
$ws.Cells.Item(45, 3).Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = "0.3907"
$ws.Cells.Item(45, 5).Value = "  +4.94%  "
$ws.Cells.Item(46, 2).Value = "Aptos"
$ws.Cells.Item(46, 3).Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = "6.923"
$ws.Cells.Item(46, 5).Value = "  +12.64%  "
$ws.Cells.Item(47, 2).Value = "Algorand"
$ws.Cells.Item(47, 3).Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = "0.1191"
$ws.Cells.Item(47, 5).Value = "  +6.46%  "
$ws.Cells.Item(48, 2).Value = "Cronos"
$ws.Cells.Item(48, 3).Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = "0.05324"
$ws.Cells.Item(48, 5).Value = "  +0.41%  "
$ws.Cells.Item(49, 5).Value = "  -1.58%  "
$ws.Cells.Item(50, 2).Value = "EnergySwap"
$ws.Cells.Item(50, 3).Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = "7.857"
$ws.Cells.Item(50, 5).Value = "  +3.17%  "
$ws.Cells.Item(51, 2).Value = "NEARProtocol"
$ws.Cells.Item(51, 3).Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = "1.271"
$ws.Cells.Item(51, 5).Value = "  +4.86%  "
